$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: update the "Date" value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- 2. Elements sheet: swap the content (and widths) of columns AK (37) and AL (38) ---
$elements = $wb.Worksheets.Item("Elements")

# Only rows whose AK/AL values actually differ need to be touched - leave
# rows where both columns already hold the same value untouched so their
# underlying cell representation (blank vs empty shared string) is not
# disturbed.
$rowsToSwap = @(1, 2, 3, 6, 8, 9, 12, 14, 15, 17, 19, 20, 23, 24)
foreach ($r in $rowsToSwap) {
    $akCell = $elements.Cells.Item($r, 37)
    $alCell = $elements.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Swap the bestFit column widths that go along with the swapped content:
# column AK (37) now holds the wide "Spécification métier..." header text,
# column AL (38) now holds the narrower "RIM Mapping" header text.
$elements.Columns.Item(37).ColumnWidth = 81.9453125
$elements.Columns.Item(38).ColumnWidth = 24.98046875
